$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9737015962600708
$ws.Range("B1").Value = 7.093147754669189
$ws.Range("C1").Value = 3.058813810348511
$ws.Range("D1").Value = 1.981501936912537
$ws.Range("E1").Value = 1.772764325141907
